# Adds the "Title and Department" project-tracker table to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Cells.Item(1, 1).Value = "ProjectName"
$ws.Cells.Item(1, 2).Value = "ClientPartner"
$ws.Cells.Item(1, 3).Value = "TeamLeads"
$ws.Cells.Item(1, 4).Value = "Status"
$ws.Cells.Item(1, 5).Value = "StartDate"
$ws.Cells.Item(1, 6).Value = "EndDate"

# Data row (row 2) - filled B:F first, then A, to match the authored
# shared-string ordering (Title added last).
$ws.Cells.Item(2, 2).Value = "Captain Planet"
$ws.Cells.Item(2, 3).Value = "The Gabosh, Captain Planet"
$ws.Cells.Item(2, 4).Value = "Active"
$ws.Cells.Item(2, 5).Value = "'2015-01-20"
$ws.Cells.Item(2, 6).Value = "'2015-01-27"
$ws.Cells.Item(2, 1).Value = "Watching Netflix Again"

# Column widths (best-fit-like) matching the authored layout.
$ws.Columns.Item(1).ColumnWidth = 15.166666666666666
$ws.Columns.Item(2).ColumnWidth = 13.166666666666666
$ws.Columns.Item(3).ColumnWidth = 24.833333333333332
$ws.Columns.Item(4).ColumnWidth = 5.666666666666667
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.666666666666666

# Selection ends on A2, as in the authored workbook.
$ws.Range("A2").Select()
